$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header cells: drop the "(*)" required-marker suffix for
# Unit/Unit Price/Currency/VAT columns (D1:G1)
$ws.Range("D1").Value = "Unit"
$ws.Range("E1").Value = "Unit Price"
$ws.Range("F1").Value = "Currency"
$ws.Range("G1").Value = "VAT"

# Clear the sample Unit/Unit Price/Currency/VAT values on the first
# data row (row 2) - the CRCC-108 code and its related values are removed
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()

# Update the view: scroll so column C is the left-most visible column
# and select H2 instead of H10
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("H2").Select()
